$d = $word.ActiveDocument

# --- 1. Turn the plain-text quackit URL (end of the "resizing images" paragraph)
#        into a real hyperlink, same as the other reference lines in the doc. ---
$resizeUrl = "https://www.quackit.com/html/howto/how_to_resize_images_in_html.cfm"

$rng = $d.Content
$found = $rng.Find.Execute($resizeUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $d.Hyperlinks.Add($rng, $resizeUrl) | Out-Null
}

# --- 2. Add a new blank paragraph, then a "column - <url>" reference paragraph
#        right after the paragraph that now holds the quackit hyperlink. ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$colUrl = "https://www.w3schools.com/howto/tryit.asp?filename=tryhow_css_two_columns_flex"
$newPara.Range.Text = "column - " + $colUrl
